# header table fix; date fix; sen msg info; comment add
#
# - Adds two new header columns: H "ФИО" (Name) and I "Комментарий" (Comment).
# - Appends 14 new data rows (17-30) below the existing 16 rows of data.
# - Every cell in the sheet is stored as plain text (General format, quote
#   look-alike numbers/dates included) to match the rest of the table, so
#   numeric-looking values ("25", "548", ...) and blank cells are written
#   with a leading apostrophe to stop Excel from auto-converting them to
#   Number/blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write $Value into $Range as TEXT, matching how the rest of this
# sheet stores every value (dates, amounts, empty cells) as t="str".
# Plain numeric strings and empty strings get auto-converted by Excel to
# Number/blank on assignment, so prefix those with an apostrophe to force
# text the same way a user typing into a text-formatted cell would.
# ---------------------------------------------------------------------------
function Set-TextCell {
    param($Range, $Value)
    if ($Value -eq "" -or $Value -match '^[+-]?\d+(\.\d+)?$') {
        $Range.Value = "'" + $Value
    } else {
        $Range.Value = $Value
    }
}

# ---------------------------------------------------------------------------
# 1) Header row: new columns H ("ФИО") and I ("Комментарий").
#    Copy the formatting of the existing header cell (G1) across so the two
#    new header cells pick up the same cell style as the rest of row 1.
# ---------------------------------------------------------------------------
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "ФИО"
$ws.Range("I1").Value = "Комментарий"

# ---------------------------------------------------------------------------
# 2) Pre-format the new data rows (17:30) so they carry the same cell style
#    as the rest of the table before we fill in values. Only the columns
#    that actually get a value on each row are touched (A:G on every new
#    row, H only on rows 19-30, I only on rows 28-30) so we don't leave
#    stray empty cells behind on rows that don't use those columns.
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").Copy() | Out-Null
$ws.Range("A17:G30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G1").Copy() | Out-Null
$ws.Range("H19:H30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G1").Copy() | Out-Null
$ws.Range("I28:I30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) New data rows.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row=17; A='19.11.2022'; B='65,88'; C='Общее'; D=''; E='Материал';          F='Skysawa';     G='' }
    @{ Row=18; A='19.11.2022'; B='25,3';  C='Общее'; D=''; E='что-нибудь такое';  F='Office';      G='' }
    @{ Row=19; A='20.12.2022'; B='25';    C='Общее'; D=''; E='Материал';          F='Skysawa';     G=''; H='Rybkin Anton' }
    @{ Row=20; A='20.12.2022'; B='25';    C='Общее'; D=''; E='Жилье';             F='MCM project'; G=''; H='Rybkin Anton | s_ryb' }
    @{ Row=21; A='20.12.2022'; B='255';   C='Общее'; D=''; E='Жилье';             F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=22; A='20.12.2022'; B='2569';  C='Общее'; D=''; E='Дичь собачья';      F='Office';      G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=23; A='20.12.2022'; B='256';   C='Общее'; D=''; E='Жилье';             F='GIPS Karpacz';G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=24; A='20.12.2022'; B='214';   C='Общее'; D=''; E='Топливо';           F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=25; A='20.12.2022'; B='548';   C='Общее'; D=''; E='Жилье';             F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=26; A='20.12.2022'; B='548';   C='Общее'; D=''; E='Жилье';             F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=27; A='20.12.2022'; B='548';   C='Общее'; D=''; E='Расходники';        F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235' }
    @{ Row=28; A='20.12.2022'; B='548';   C='Общее'; D=''; E='Расходники';        F='MCM project'; G=''; H='Rybkin Anton; s_ryb; 192713235'; I='ну типа того вот так и эдак' }
    @{ Row=29; A='20.12.2022'; B='4123';  C='Общее'; D=''; E='Топливо';           F='Karpacz';     G=''; H='Rybkin Anton; s_ryb';            I='шла саша по шосе и сосала сушку' }
    @{ Row=30; A='20.12.2022'; B='589';   C='Общее'; D=''; E='херня какая-то';    F='GIPS Karpacz';G=''; H='Rybkin Anton; s_ryb';            I='нужная херня очень' }
)

foreach ($r in $rows) {
    Set-TextCell $ws.Range("A" + $r.Row) $r.A
    Set-TextCell $ws.Range("B" + $r.Row) $r.B
    Set-TextCell $ws.Range("C" + $r.Row) $r.C
    Set-TextCell $ws.Range("D" + $r.Row) $r.D
    Set-TextCell $ws.Range("E" + $r.Row) $r.E
    Set-TextCell $ws.Range("F" + $r.Row) $r.F
    Set-TextCell $ws.Range("G" + $r.Row) $r.G
    if ($r.ContainsKey('H')) {
        Set-TextCell $ws.Range("H" + $r.Row) $r.H
    }
    if ($r.ContainsKey('I')) {
        Set-TextCell $ws.Range("I" + $r.Row) $r.I
    }
}
